$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 2
    $ws.Range("C2").Value = "丽水·新年动漫狂欢盛典"
    $ws.Range("F2").Value = 278
    $ws.Range("G2").Value = 4500

    # Row 3
    $ws.Range("C3").Value = "龙泉·崩X铁X原ONLY"
    $ws.Range("G3").Value = 5000

    # Row 4
    $ws.Range("C4").Value = "丽水·YA●怀旧only"
    $ws.Range("G4").Value = 3500

    # Row 5
    $ws.Range("C5").Value = "丽水·LPJ 现实X次元动漫展"
    $ws.Range("G5").Value = 4500
}
